$wb = $excel.ActiveWorkbook

# --- Rename the review sheet (tab) from REVIEW-SHEET to LH_WF_NAVIGATION_REVIEWS ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LH_WF_NAVIGATION_REVIEWS"

# --- All reviewer verifications are now closed ---
$ws1.Range("I2").Value = "closed"
$ws1.Range("I3").Value = "closed"
$ws1.Range("I4").Value = "closed"

# --- Log the change in the VERSION-HISTORY sheet as a new row (v1.2) ---
$ws2 = $wb.Worksheets.Item(2)

# Copy formatting from the row above so the new row matches the existing table style
[void]$ws2.Range("A3:D3").Copy()
[void]$ws2.Range("A4:D4").PasteSpecial(-4122)

$ws2.Range("A4").Value = "v1.2"
$ws2.Range("B4").Value = "Hala Eldaly"
$ws2.Range("C4").Value = "All reviewer verifications are closed, and the file tab is renamed"
$ws2.Range("D4").Value = $ws2.Range("D3").Value2

# --- Restore cursor/selection state on each sheet (as left by the author) ---
[void]$ws1.Select()
[void]$ws1.Range("E11").Select()
[void]$ws2.Select()
[void]$ws2.Range("C14").Select()
